$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Existing title cell text changed from "Contract" to "Contract conditions" ---
$ws.Range("A2").Value = "Contract conditions"

# --- Re-point existing contractor rows onto their (unchanged) formulas ---
$ws.Range("A4").Value = "Contractor"
$ws.Range("B5").Value = "Name"
$ws.Range("C5").Value = '${ctx[''contract''].contractor.name}'
$ws.Range("B6").Value = "Birth date"
$ws.Range("C6").Value = '${ctx[''contract''].contractor.birthDate}'

$ws.Range("A8").Value = "Beneficiaries"
$ws.Range("B9").Value = '${beneficiary.name}'
$ws.Range("C9").Value = '$beneficiary.phoneNumber}'

# --- New "Fees" block (columns K:O, rows 4-7) ---
$ws.Range("K4").Value = "Fees"
$ws.Range("K4").Font.Bold = $true

$ws.Range("K5").Value = "Non-benefiary calls"
$ws.Range("L5").Value = "Weekdays"
$ws.Range("M5").Value = "08:00-15:59"
$ws.Range("N5").Value = 30
$ws.Range("O5").Value = "HUF"

$ws.Range("K6").Value = "Non-benefiary calls"
$ws.Range("L6").Value = "Weekdays"
$ws.Range("M6").Value = "16:00-07:59"
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = "HUF"

$ws.Range("K7").Value = "Benefiary calls"
$ws.Range("L7").Value = "Everyday"
$ws.Range("M7").Value = "00:00-23:59"
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = "HUF"

# --- Column widths for the new columns (engine stores width = ColumnWidth + 5/6) ---
$ws.Range("K:K").ColumnWidth = 16.053385416666668
$ws.Range("L:M").ColumnWidth = 9.944010416666666

# --- Selection moves to K7 ---
$null = $ws.Range("K7").Select()

# --- Comment text / author on B9 ---
$excel.UserName = 'Szerző'
$comment = $ws.Range("B9").Comment
$null = $comment.Text('Szerző:' + "`n" + 'jx:each(items="ctx[''contract''].beneficiaries", groupBy="beneficiary.name", var="beneficiary", lastCell="C9")')
try { $comment.Author = 'Szerző' } catch { }

# --- Page setup for the (previously blank) Munka2 / Munka3 sheets ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

$ws3 = $wb.Worksheets.Item(3)
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1
